# Update cryptocurrency market data (refresh dated 2023-05-24)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 26354
$ws.Range("E2").Value = 510900629657
$ws.Range("F2").Value = 17010212468
$ws.Range("G2").Value = -3.09514

$ws.Range("D3").Value = 1800.15
$ws.Range("E3").Value = 216516951385
$ws.Range("F3").Value = 8056153470
$ws.Range("G3").Value = -2.85654

$ws.Range("D4").Value = 0.999027
$ws.Range("E4").Value = 82946120777
$ws.Range("F4").Value = 19637127698
$ws.Range("G4").Value = -0.07069

$ws.Range("D5").Value = 305.82
$ws.Range("E5").Value = 48291483308
$ws.Range("F5").Value = 566920174
$ws.Range("G5").Value = -2.41587

$ws.Range("D6").Value = 0.9994960000000001
$ws.Range("E6").Value = 29128430066
$ws.Range("F6").Value = 4136397359
$ws.Range("G6").Value = 0.00287

$ws.Range("D7").Value = 0.454593
$ws.Range("E7").Value = 23586402042
$ws.Range("F7").Value = 1118162463
$ws.Range("G7").Value = -2.26386

$ws.Range("D8").Value = 0.364502
$ws.Range("E8").Value = 12773672464
$ws.Range("F8").Value = 189220373
$ws.Range("G8").Value = -1.64733

$ws.Range("D9").Value = 1799.66
$ws.Range("E9").Value = 12003448925
$ws.Range("F9").Value = 11872853
$ws.Range("G9").Value = -2.73351

$ws.Range("D10").Value = 0.07085
$ws.Range("E10").Value = 9885372739
$ws.Range("F10").Value = 318146008
$ws.Range("G10").Value = -2.80767

$ws.Range("D11").Value = 0.8743649999999999
$ws.Range("E11").Value = 8115102468
$ws.Range("F11").Value = 240732378
$ws.Range("G11").Value = -1.78461

$ws.Range("D12").Value = 19.21
$ws.Range("E12").Value = 7614453735
$ws.Range("F12").Value = 320500166
$ws.Range("G12").Value = -4.05471

$ws.Range("D13").Value = 0.076942
$ws.Range("E13").Value = 6947582660
$ws.Range("F13").Value = 371089328
$ws.Range("G13").Value = -2.17915

$ws.Range("E14").Value = 6497395659
$ws.Range("F14").Value = 104956316
$ws.Range("G14").Value = -2.75013

$ws.Range("D15").Value = 86.02
$ws.Range("E15").Value = 6281708863
$ws.Range("F15").Value = 844838988
$ws.Range("G15").Value = -5.94447

$ws.Range("D16").Value = 0.998819
$ws.Range("E16").Value = 5326507249
$ws.Range("F16").Value = 2062154645
$ws.Range("G16").Value = -0.09579

$ws.Range("E17").Value = 5041769039
$ws.Range("F17").Value = 128250961
$ws.Range("G17").Value = -4.19097

$ws.Range("E18").Value = 4728612657
$ws.Range("F18").Value = 139738440
$ws.Range("G18").Value = -3.85001

$ws.Range("D19").Value = 0.999757
$ws.Range("E19").Value = 4627726121
$ws.Range("F19").Value = 120167253
$ws.Range("G19").Value = -0.01026

$ws.Range("D20").Value = 26389
$ws.Range("E20").Value = 4119525452
$ws.Range("F20").Value = 78228852
$ws.Range("G20").Value = -3.10933

$ws.Range("E21").Value = 3751244353
$ws.Range("F21").Value = 46611676
$ws.Range("G21").Value = -1.9163

$ws.Range("E22").Value = 3292517226
$ws.Range("F22").Value = 223259
$ws.Range("G22").Value = 0.76259

$ws.Range("D23").Value = 6.33
$ws.Range("E23").Value = 3274326621
$ws.Range("F23").Value = 154853396
$ws.Range("G23").Value = -2.68206

$ws.Range("E24").Value = 3045991777
$ws.Range("F24").Value = 77146006
$ws.Range("G24").Value = -1.01912

$ws.Range("D25").Value = 1.96
$ws.Range("E25").Value = 2895113694
$ws.Range("F25").Value = 12120677
$ws.Range("G25").Value = -3.38213

$ws.Range("D26").Value = 150.2
$ws.Range("E26").Value = 2725867831
$ws.Range("F26").Value = 64594444
$ws.Range("G26").Value = -0.82007

$ws.Range("D27").Value = 44.76
$ws.Range("E27").Value = 2686737474
$ws.Range("F27").Value = 7571353
$ws.Range("G27").Value = -2.74835

$ws.Range("D28").Value = 17.78
$ws.Range("E28").Value = 2509993920
$ws.Range("F28").Value = 75979091
$ws.Range("G28").Value = -3.26488

$ws.Range("D29").Value = 0.08663
$ws.Range("E29").Value = 2322029725
$ws.Range("F29").Value = 37245886
$ws.Range("G29").Value = -1.81561

$ws.Range("D30").Value = 112.09
$ws.Range("E30").Value = 2175646447
$ws.Range("F30").Value = 64931074
$ws.Range("G30").Value = -2.93662

$ws.Range("E31").Value = 2102657365
$ws.Range("F31").Value = 27694350
$ws.Range("G31").Value = -4.33569

$ws.Range("D32").Value = 0.998657
$ws.Range("E32").Value = 2039285129
$ws.Range("F32").Value = 231734314
$ws.Range("G32").Value = -0.10502

$ws.Range("E33").Value = 1905139299
$ws.Range("F33").Value = 110167715
$ws.Range("G33").Value = -1.06639

$ws.Range("E34").Value = 1806177261
$ws.Range("F34").Value = 64111860
$ws.Range("G34").Value = 0.418

$ws.Range("D35").Value = 0.051021
$ws.Range("E35").Value = 1605362166
$ws.Range("F35").Value = 20034280
$ws.Range("G35").Value = -2.73254

$ws.Range("E36").Value = 1601229647
$ws.Range("F36").Value = 67298555
$ws.Range("G36").Value = -4.84701

$ws.Range("D37").Value = 0.059954
$ws.Range("E37").Value = 1514749538
$ws.Range("F37").Value = 7729901
$ws.Range("G37").Value = -3.52959

$ws.Range("D38").Value = 100.64
$ws.Range("E38").Value = 1463902987
$ws.Range("F38").Value = 16058149
$ws.Range("G38").Value = -2.05009

$ws.Range("E39").Value = 1433723896
$ws.Range("F39").Value = 61619698
$ws.Range("G39").Value = -3.955

$ws.Range("E40").Value = 1418462287
$ws.Range("F40").Value = 191198559
$ws.Range("G40").Value = -4.69035

$ws.Range("D41").Value = 0.01944619
$ws.Range("E41").Value = 1414436619
$ws.Range("F41").Value = 41258243
$ws.Range("G41").Value = 0.19772

$ws.Range("E42").Value = 1229059095
$ws.Range("F42").Value = 73263611
$ws.Range("G42").Value = -4.77539

$ws.Range("D43").Value = 0.155172
$ws.Range("E43").Value = 1124299133
$ws.Range("F43").Value = 48688409
$ws.Range("G43").Value = -5.01463

$ws.Range("D44").Value = 0.11669
$ws.Range("E44").Value = 1049861239
$ws.Range("F44").Value = 38921261
$ws.Range("G44").Value = -4.54035

$ws.Range("D45").Value = 0.09041299999999999
$ws.Range("E45").Value = 1035557832
$ws.Range("F45").Value = 41873
$ws.Range("G45").Value = -0.10993

$ws.Range("D46").Value = 0.999812
$ws.Range("E46").Value = 1018906889
$ws.Range("F46").Value = 22087900
$ws.Range("G46").Value = 0.11962

$ws.Range("D47").Value = 0.9976390000000001
$ws.Range("E47").Value = 1002267460
$ws.Range("F47").Value = 7708758
$ws.Range("G47").Value = -0.13573

$ws.Range("D48").Value = 0.999177
$ws.Range("E48").Value = 999122698
$ws.Range("F48").Value = 5645718
$ws.Range("G48").Value = -0.09465

$ws.Range("D49").Value = 2.67
$ws.Range("E49").Value = 981014586
$ws.Range("F49").Value = 202300461
$ws.Range("G49").Value = -0.42676

$ws.Range("D50").Value = 0.839825
$ws.Range("E50").Value = 930208452
$ws.Range("F50").Value = 97593482
$ws.Range("G50").Value = -2.5014

$ws.Range("D51").Value = 0.330898
$ws.Range("E51").Value = 923775935
$ws.Range("F51").Value = 273584511
$ws.Range("G51").Value = -9.167120000000001

